$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 76; this shifts the existing rows 76-106 down
# to become rows 77-107 (dimension grows from A1:R106 to A1:R107).
$ws.Rows("76:76").Insert()

# Populate the newly-inserted row 76 with the new weekly price record.
$ws.Cells.Item(76, 1).Value = 7
$ws.Cells.Item(76, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(76, 3).Value = "Ñuble"
$ws.Cells.Item(76, 4).Value = 44875
$ws.Cells.Item(76, 5).Value = 16
$ws.Cells.Item(76, 6).Value = 100112021
$ws.Cells.Item(76, 7).Value = "Ají"
$ws.Cells.Item(76, 8).Value = "Inferno"
$ws.Cells.Item(76, 9).Value = "Primera"
$ws.Cells.Item(76, 10).Value = 50
$ws.Cells.Item(76, 11).Value = 20000
$ws.Cells.Item(76, 12).Value = 20000
$ws.Cells.Item(76, 13).Value = 20000
$ws.Cells.Item(76, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(76, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(76, 16).Value = 2000
$ws.Cells.Item(76, 17).Value = 10
$ws.Cells.Item(76, 18).Value = "Hortaliza"
